# Insert a new weekly price record at row 98 for
# "Feria Lagunitas de Puerto Montt - Brócoli" (Hortaliza), pushing every
# existing record from row 98 down by one row (98->99, ..., 182->183).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 98:182 down to 99:183, copying formatting from the row above
# (this also grows the sheet's used range / dimension to A1:R183).
$ws.Rows.Item(98).Insert()

# Populate the newly-inserted row 98 with the new record's data.
$ws.Cells.Item(98, 1).Value = 4
$ws.Cells.Item(98, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(98, 3).Value = "Los Lagos"
$ws.Cells.Item(98, 4).Value = 44447
$ws.Cells.Item(98, 5).Value = 10
$ws.Cells.Item(98, 6).Value = 100112023
$ws.Cells.Item(98, 7).Value = "Brócoli"
$ws.Cells.Item(98, 8).Value = "Sin especificar"
$ws.Cells.Item(98, 9).Value = "Segunda"
$ws.Cells.Item(98, 10).Value = 100
$ws.Cells.Item(98, 11).Value = 1000
$ws.Cells.Item(98, 12).Value = 1000
$ws.Cells.Item(98, 13).Value = 1000
$ws.Cells.Item(98, 14).Value = "$/unidad"
$ws.Cells.Item(98, 15).Value = "Región Metropolitana"
$ws.Cells.Item(98, 16).Value = 1000
$ws.Cells.Item(98, 17).Value = 1
$ws.Cells.Item(98, 18).Value = "Hortaliza"
